$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.975.39"
$ws.Range("E2").Value = "  +2.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.702.14"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.20"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3970"
$ws.Range("E7").Value = "  +1.89%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4024"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -1.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.82"
$ws.Range("E10").Value = "  +0.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08799"
$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.90"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.451"
$ws.Range("E14").Value = "  -0.46%  "

$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001349"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.708.57"
$ws.Range("E17").Value = "  +1.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.28"
$ws.Range("E18").Value = "  -1.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07184"
$ws.Range("E19").Value = "  -1.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.59"
$ws.Range("E20").Value = "  +3.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.345"
$ws.Range("E21").Value = "  +1.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.36"
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.978.07"
$ws.Range("E24").Value = "  +2.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.346"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.935"
$ws.Range("E26").Value = "  -3.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.71"
$ws.Range("E27").Value = "  +5.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.154"
$ws.Range("E28").Value = "  +14.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.60"
$ws.Range("E29").Value = "  -3.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "149.60"
$ws.Range("E30").Value = "  +8.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.307"
$ws.Range("E31").Value = "  -4.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.624"
$ws.Range("E32").Value = "  +24.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.900.50"
$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08541"
$ws.Range("E34").Value = "  -2.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03130"
$ws.Range("E35").Value = "  +3.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.038"
$ws.Range("E36").Value = "  -0.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.113"
$ws.Range("E37").Value = "  -3.05%  "

$ws.Range("E38").Value = "  +2.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.87"
$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09526"
$ws.Range("E40").Value = "  +4.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8213"
$ws.Range("E41").Value = "  +2.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.93"
$ws.Range("E42").Value = "  -1.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.476"
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.10"
$ws.Range("E44").Value = "  -2.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.668"
$ws.Range("E45").Value = "  +1.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7370"
$ws.Range("E46").Value = "  +1.82%  "

$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.397"
$ws.Range("E48").Value = "  -1.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08738"
$ws.Range("E49").Value = "  +8.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.003"
$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.03"
$ws.Range("E51").Value = "  +0.01%  "
